# Apply "add single and multi corrector" edit:
# Shifts data left by one column (drops old column A "NO"/serial-number
# prefixes), splits merged/garbled strings back into their proper
# Item / Result / Unit / Reference-range / Method columns, and fixes a
# couple of values/typos along the way. Final used range becomes A1:E6
# (previously A1:F6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "项目"
$ws.Range("B1").Value = "结果"
$ws.Range("C1").Value = "单位"
$ws.Range("D1").Value = "参考区间"
$ws.Range("E1").Value = "实验方法"

# --- Row 2: 乙肝表面抗原 ---
$ws.Range("A2").Value = "乙肝表面抗原"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "（"
$ws.Range("D2").Value = "0.00-0.05"
$ws.Range("E2").Value = "化学发光法"

# --- Row 3: 抗-HBe抗体 ---
$ws.Range("A3").Value = "抗-HBe抗体"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.29"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "0.00-10.00"
$ws.Range("E3").Value = "化学发光法"

# --- Row 4: 乙肝病毒e抗原 ---
$ws.Range("A4").Value = "乙肝病毒e抗原"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "s/co"
$ws.Range("D4").Value = "0.00-1.00"
$ws.Range("E4").Value = "化学发光法"

# --- Row 5: 抗-HBe抗体 ---
$ws.Range("A5").Value = "抗-HBe抗体"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "67.52"
$ws.Range("C5").Value = "s/co"
$ws.Range("D5").Value = "1.00-999.00"
$ws.Range("E5").Value = "化学发光法"

# --- Row 6: 抗-HBe抗体 ---
$ws.Range("A6").Value = "抗-HBe抗体"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "4.19"
$ws.Range("C6").Value = "S/co"
$ws.Range("D6").Value = "0.00-1.00"
$ws.Range("E6").Value = "化学发光法"

# --- Clear the now-unused column F (dimension shrinks from A1:F6 to A1:E6) ---
$ws.Range("F1:F6").Clear()
